# Update chapter {chapter} - {title}
#
# Slide 4 ("Why use new terminology?") has a bullet reading:
#   db.help() or db.stats()
# The two runs that make up "db.stats()" are collapsed into a single
# run containing just a tab character (formatted as Courier New, Italian
# input language), matching the authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$fullText = $tr.Text
$target = "db.stats()"
$pos = $fullText.IndexOf($target)

if ($pos -ge 0) {
    $startChar = $pos + 1
    $len = $target.Length
    $rng = $tr.Characters($startChar, $len)

    $rng.Text = [char]9
}
